# Added scraper for NF group for St. Petersburg:
# append a new row (row 9) to the Sales sheet, re-using the same listing
# (id 10056376) re-scraped on a later date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate row 2 (same listing) into the new row 9, including blank cells,
# then update the scrape_date and sale_date columns to reflect the new scrape.
$ws.Range("A2:V2").Copy($ws.Range("A9:V9"))

$ws.Cells.Item(9, 13).Value = "23_10_2023"   # scrape_date
$ws.Cells.Item(9, 22).Value = "25/10/2023"   # sale_date
